# Update cryptocurrency price/volume figures to match the refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (kept as text so values such as
# "275.00" or "-2.12%" are not silently reinterpreted as numbers/percentages).
$updates = [ordered]@{
    "D2" = "275.00"
    "E2" = "-2.12%"
    "D3" = "27.20"
    "E3" = "1.30%"
    "D4" = "4.752"
    "E4" = "-3.87%"
    "D5" = "0.06303"
    "E5" = "-1.68%"
    "D6" = "6.925"
    "E6" = "-0.89%"
    "D7" = "1.352"
    "E7" = "33.09%"
    "D8" = "0.8767"
    "E8" = "-1.16%"
    "D9" = "0.1510"
    "E9" = "1.43%"
    "D10" = "0.05082"
    "E10" = "-2.17%"
    "D11" = "0.07599"
    "E11" = "2.38%"
    "D12" = "0.02978"
    "E12" = "-4.22%"
    "D13" = "0.09010"
    "E13" = "-0.42%"
    "D14" = "0.001571"
    "E14" = "0.72%"
    "D15" = "0.0006342"
    "E15" = "0.14%"
    "D16" = "0.006002"
    "E16" = "-0.40%"
    "D17" = "3.446"
    "E17" = "-1.39%"
    "D18" = "3.298"
    "E18" = "-1.64%"
    "D19" = "2.285"
    "E19" = "-0.61%"
    "D21" = "0.1336"
    "E21" = "0.43%"
    "D22" = "3.916"
    "E22" = "-0.42%"
    "D23" = "0.04399"
    "E23" = "1.05%"
    "D24" = "0.001173"
    "E24" = "-0.61%"
    "E25" = "4.23%"
    "D26" = "0.0001200"
    "E26" = "-0.06%"
    "D27" = "0.0001937"
    "E27" = "14.32%"
    "D40" = "0.04099"
    "E40" = "0.00%"
    "D41" = "0.006791"
    "E41" = "2.05%"
    "D42" = "0.1173"
    "E42" = "-0.40%"
    "D43" = "0.002129"
    "E43" = "-9.80%"
    "D44" = "0.01153"
    "E44" = "-11.75%"
    "D45" = "0.00005165"
    "E45" = "-1.49%"
    "D47" = "0.02300"
    "E47" = "2.23%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

